$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Fill in grade values for row 12 (student "Hhhhhhh Stanisław")
$ws.Range("C12").Value = 2.5
$ws.Range("D12").Value = 4.5
$ws.Range("H12").Value = 5
$ws.Range("J12").Value = 4.5
$ws.Range("L12").Value = 4.5
$ws.Range("N12").Value = 4.5
$ws.Range("O12").Value = 10

# Update the view state: frozen-pane top-left cell and active selection
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("G12").Select()
